$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 999.3333
$ws.Range("J9").Value = 999.3333
$ws.Range("L9").Value = 999.3333
$ws.Range("N9").Value = -1337.3333
$ws.Range("H11").Value = 68
$ws.Range("I11").Value = 68
$ws.Range("K11").Value = 68
$ws.Range("M11").Value = 72
$ws.Range("H42").Value = 5
$ws.Range("I42").Value = 5
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 15
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 215
$ws.Range("N42").ClearContents()
$ws.Range("H64").Value = 3790.3
$ws.Range("I64").Value = 3790.3
$ws.Range("K64").Value = 3790.3
$ws.Range("M64").Value = -3542.3
$ws.Range("H67").Value = 3790.3
$ws.Range("I67").Value = 3790.3
$ws.Range("K67").Value = 3790.3
$ws.Range("M67").Value = -2932.3
$ws.Range("H76").Value = 3003
$ws.Range("I76").Value = 3003
$ws.Range("K76").Value = 3003
$ws.Range("M76").Value = -2688
$ws.Range("H79").Value = 3003
$ws.Range("I79").Value = 3003
$ws.Range("K79").Value = 3003
$ws.Range("M79").Value = -1911
$ws.Range("H105").Value = 20000
$ws.Range("J105").Value = 20000
$ws.Range("L105").Value = 20000
$ws.Range("N105").Value = -26988
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8172.3335
$ws.Range("H45").Value = 3099.8572
$ws.Range("J45").Value = 3860
$ws.Range("L45").Value = 3860
$ws.Range("N45").Value = -4614
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 1922.2222
$ws.Range("I132").Value = 1850
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5550
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -3020
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 654.55554
$ws.Range("I22").Value = 486.375
$ws.Range("K22").Value = 486.375
$ws.Range("M22").Value = -313.375
$ws.Range("H94").Value = 2054.5625
$ws.Range("I94").Value = 3393.4285
$ws.Range("J94").Value = 1013.2222
$ws.Range("K94").Value = 3393.4285
$ws.Range("L94").Value = 1013.2222
$ws.Range("M94").Value = -2942.4285
$ws.Range("N94").Value = -1915.2222
$ws.Range("H99").Value = 2190.7222
$ws.Range("J99").Value = 3451.75
$ws.Range("L99").Value = 3451.75
$ws.Range("N99").Value = -6447.75
$ws.Range("H134").Value = 1107.0769
$ws.Range("I134").Value = 1190.75
$ws.Range("J134").Value = 973.2
$ws.Range("K134").Value = 3572.25
$ws.Range("L134").Value = 2919.6
$ws.Range("M134").Value = -1037.25
$ws.Range("N134").Value = -7989.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 17762.143
$ws.Range("J4").Value = 17762.143
$ws.Range("L4").Value = 17762.143
$ws.Range("N4").Value = -17986.143
$ws.Range("H86").Value = 224983.6
$ws.Range("I86").Value = 224983.6
$ws.Range("K86").Value = 224983.6
$ws.Range("M86").Value = -223860.6
$ws.Range("H89").Value = 224983.6
$ws.Range("I89").Value = 224983.6
$ws.Range("K89").Value = 1124918
$ws.Range("M89").Value = -1119302
$ws.Range("H132").Value = 2998
$ws.Range("I132").Value = 1999.5
$ws.Range("K132").Value = 5998.5
$ws.Range("M132").Value = -3468.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 400000350
$ws.Range("I4").Value = 400000350
$ws.Range("K4").Value = 1200001050
$ws.Range("M4").Value = -1200000938
$ws.Range("H9").Value = 2324.5
$ws.Range("J9").Value = 2324.5
$ws.Range("L9").Value = 6973.5
$ws.Range("N9").Value = -7421.5
$ws.Range("H102").Value = 11250
$ws.Range("J102").Value = 11250
$ws.Range("L102").Value = 33750
$ws.Range("N102").Value = -38618
$ws.Range("H103").Value = 32468.625
$ws.Range("I103").Value = 84100
$ws.Range("J103").Value = 1489.8
$ws.Range("K103").Value = 252300
$ws.Range("L103").Value = 4469.4
$ws.Range("M103").Value = -251421
$ws.Range("N103").Value = -6227.4
$ws.Range("H107").Value = 1362
$ws.Range("J107").Value = 1362
$ws.Range("L107").Value = 4086
$ws.Range("N107").Value = -7926
$ws.Range("H134").Value = 1212.25
$ws.Range("I134").Value = 1212.25
$ws.Range("K134").Value = 3636.75
$ws.Range("M134").Value = 1433.25
$ws.Range("H139").Value = 5566.6
$ws.Range("I139").Value = 4458.25
$ws.Range("K139").Value = 13374.75
$ws.Range("M139").Value = -8234.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3802.7334
$ws.Range("I97").Value = 3808
$ws.Range("J97").Value = 3792.2
$ws.Range("K97").Value = 3808
$ws.Range("L97").Value = 3792.2
$ws.Range("M97").Value = -3312
$ws.Range("N97").Value = -4784.2
$ws.Range("H102").Value = 1391.2
$ws.Range("J102").Value = 1200
$ws.Range("L102").Value = 1200
$ws.Range("N102").Value = -4444
$ws.Range("H126").Value = 8753
$ws.Range("I126").Value = 6670.6665
$ws.Range("K126").Value = 20011.9995
$ws.Range("M126").Value = -17541.9995
$ws.Range("H132").Value = 3337.1765
$ws.Range("I132").Value = 3295.75
$ws.Range("K132").Value = 9887.25
$ws.Range("M132").Value = -7357.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1558.7142
$ws.Range("I40").Value = 1558.7142
$ws.Range("K40").Value = 1558.7142
$ws.Range("M40").Value = -1422.7142
$ws.Range("H61").Value = 2831.3333
$ws.Range("I61").Value = 2997.8
$ws.Range("J61").Value = 1999
$ws.Range("K61").Value = 2997.8
$ws.Range("L61").Value = 1999
$ws.Range("M61").Value = -2795.8
$ws.Range("N61").Value = -2403
$ws.Range("H113").Value = 2831.3333
$ws.Range("I113").Value = 2997.8
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 2997.8
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = -827.8000000000002
$ws.Range("N113").Value = -6339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3169366.2
$ws.Range("I100").Value = 5362267.5
$ws.Range("J100").Value = 1842.3334
$ws.Range("K100").Value = 10724535
$ws.Range("L100").Value = 3684.6668
$ws.Range("M100").Value = -10723994
$ws.Range("N100").Value = -4766.6668
$ws.Range("H126").Value = 1096.6666
$ws.Range("I126").Value = 1096.6666
$ws.Range("K126").Value = 3289.9998
$ws.Range("M126").Value = -819.9998000000001

Write-Host "Applied all changes"
